$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Raw value changes (row 12) ---
$ws.Range("G12").Value2 = 1240524717.0500007
$ws.Range("I12").Value2 = 4188377156

# --- Raw value changes (row 13) ---
$ws.Range("G13").Value2 = 319819483.18000001
$ws.Range("I13").Value2 = 1012006300

# --- Raw value changes (row 14) ---
$ws.Range("G14").Value2 = 34063116.800000042
$ws.Range("I14").Value2 = -4419159.29

# --- Raw value changes (row 16) ---
$ws.Range("G16").Value2 = -60834434.380000003
$ws.Range("I16").Value2 = -162861893.59999999

# --- Row 18: G18 becomes a formula (SUM), I18 recalculates automatically ---
$ws.Range("G18").Formula = "=SUM(G12:G17)"

# --- Raw value changes (row 19) ---
$ws.Range("G19").Value2 = -379300000.00000012
$ws.Range("I19").Value2 = -1160500000

# --- Row 21: G21 becomes a formula (SUM), I21 recalculates automatically ---
$ws.Range("G21").Formula = "=SUM(G18:G20)"

# --- Raw value changes (row 22) ---
$ws.Range("G22").Value2 = -20015625

# --- Raw value changes (row 26) ---
$ws.Range("G26").Value2 = 1029174575.116062
$ws.Range("I26").Value2 = 1010658959

$excel.CalculateFull()
